$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells B5, B6, B7, B8, B11, B12, B18 go from the "?%" text placeholder
# to an actual numeric 90% value, formatted as a percentage (matches
# the existing style used elsewhere in column B, numFmtId 9 -> "0%").
$percentCells = @("B5", "B6", "B7", "B8", "B11", "B12", "B18")
foreach ($addr in $percentCells) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "0%"
    $rng.Value = 0.9
}

# Add the two new note rows at the bottom of the sheet.
$ws.Range("A33").Value = "Design Pattern Applied: Factory, Object Pool, Type Object, Observer"
$ws.Range("A34").Value = "All Apply :  Factory, Object Pool, Type Object, Observer, State, Singleton"

# Update the active selection to reflect where the author left off editing.
$ws.Range("D26").Select()
